# Updated symbol list on Wed Dec 21 06:56:30 UTC 2022 with GitHub Actions
#
# Applies the latest coinranking.com price/volume refresh to the "Price"
# (column D) and "Volume(1h)" (column E) columns. Values are stored as
# plain text in the source sheet (not numbers), so each cell's number
# format is forced to Text ("@") before the write and then cleared again
# immediately afterwards so no new/lingering cell style is left behind -
# this keeps the write a pure value change, matching how the sheet looked
# before the refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2";  Value = "248.85" },
    @{ Cell = "D3";  Value = "22.58" },
    @{ Cell = "D4";  Value = "5.416" },
    @{ Cell = "D7";  Value = "6.321" },
    @{ Cell = "D8";  Value = "0.8124" },
    @{ Cell = "D9";  Value = "0.9253" },
    @{ Cell = "D10"; Value = "0.1424" },
    @{ Cell = "D11"; Value = "0.07416" },
    @{ Cell = "D12"; Value = "0.03023" },
    @{ Cell = "D13"; Value = "0.03019" },
    @{ Cell = "D15"; Value = "3.760" },
    @{ Cell = "D16"; Value = "0.001579" },
    @{ Cell = "D17"; Value = "0.04756" },
    @{ Cell = "D19"; Value = "0.0005791" },
    @{ Cell = "E19"; Value = "18OneONE" },
    @{ Cell = "D20"; Value = "0.006446" },
    @{ Cell = "D21"; Value = "0.004997" },
    @{ Cell = "D22"; Value = "0.001025" },
    @{ Cell = "D27"; Value = "0.1302" },
    @{ Cell = "D40"; Value = "0.03996" },
    @{ Cell = "D41"; Value = "0.006821" },
    @{ Cell = "D42"; Value = "0.1066" },
    @{ Cell = "D43"; Value = "0.002711" },
    @{ Cell = "D44"; Value = "0.007508" },
    @{ Cell = "D45"; Value = "0.00005923" },
    @{ Cell = "D47"; Value = "0.4301" },
    @{ Cell = "E47"; Value = "46CoinbaseStockTokenCOINWorstin24h" },
    @{ Cell = "D48"; Value = "0.2107" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Force Text format so a numeric-looking string (e.g. "248.85") is
    # written back as text instead of being auto-converted to a number.
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    # Drop the temporary "@" format again so the cell's style index is left
    # exactly as it was (General / default style) - only the value changes.
    $rng.ClearFormats()
}
